$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 108-114 append the September 2021 auction results below the
# existing data (rows 2-107). Column A holds the auction date (stored as
# plain text, matching the existing "dd-mm-yyyy" shared strings already in
# the sheet), columns B-G hold the numeric figures.
#
# Note: dates whose leading day-of-month component is <= 12 (08, 09, 10)
# are ambiguous with a month number, so a bare string assignment would be
# auto-converted into a real date serial by the host. Pre-formatting those
# three cells as Text ("@") before the assignment keeps them as literal
# strings; the format is reset back to Normal immediately afterwards so the
# cells end up with no explicit style, exactly like all of the other date
# cells in column A.
$ws.Range("A108:A110").NumberFormat = "@"

$ws.Range("A108").Value = "08-09-2021"
$ws.Range("B108").Value = 600000
$ws.Range("C108").Value = 908000
$ws.Range("D108").Value = 900000
$ws.Range("E108").Value = 692000
$ws.Range("F108").Value = 208000
$ws.Range("G108").Value = 1.49

$ws.Range("A109").Value = "09-09-2021"
$ws.Range("B109").Value = 800000
$ws.Range("C109").Value = 1485000
$ws.Range("D109").Value = 800000
$ws.Range("E109").Value = 535000
$ws.Range("F109").Value = 265000
$ws.Range("G109").Value = 1.47

$ws.Range("A110").Value = "10-09-2021"
$ws.Range("B110").Value = 800000
$ws.Range("C110").Value = 910000
$ws.Range("D110").Value = 800000
$ws.Range("E110").Value = 490000
$ws.Range("F110").Value = 310000
$ws.Range("G110").Value = 1.49

$ws.Range("A108:A110").Style = "Normal"

$ws.Range("A111").Value = "13-09-2021"
$ws.Range("B111").Value = 300000
$ws.Range("C111").Value = 770000
$ws.Range("D111").Value = 450000
$ws.Range("E111").Value = 305000
$ws.Range("F111").Value = 145000
$ws.Range("G111").Value = 1.47

$ws.Range("A112").Value = "14-09-2021"
$ws.Range("B112").Value = 300000
$ws.Range("C112").Value = 629000
$ws.Range("D112").Value = 450000
$ws.Range("E112").Value = 380000
$ws.Range("F112").Value = 70000
$ws.Range("G112").Value = 1.43

$ws.Range("A113").Value = "15-09-2021"
$ws.Range("B113").Value = 800000
$ws.Range("C113").Value = 665000
$ws.Range("D113").Value = 400000
$ws.Range("E113").Value = 335000
$ws.Range("F113").Value = 65000
$ws.Range("G113").Value = 1.48

$ws.Range("A114").Value = "16-09-2021"
$ws.Range("B114").Value = 800000
$ws.Range("D114").Value = 0
